# Updates cryptos list values/percentages per latest scrape; ImmutableX/NEARProtocol rows swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.923.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.909.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.73%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.32"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.83%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.506"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.908.86"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.78"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.144"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.16"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.387.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.779.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.904.47"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.672"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.11"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.66"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.93"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.52%  "

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.22"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.57"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.86%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0839"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.68"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.04"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.32"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.123"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.73"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.292"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.19"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0347"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "373.77"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.94"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.659.14"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.41%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.05"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.30%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.19%  "
